$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 208 (Excel shifts rows 208-220 down to 209-221)
$ws.Rows("208:208").Insert()

# Populate the new row 208 with the new weekly data record
$ws.Range("A208").Value = 10
$ws.Range("B208").Value = "Vega Modelo de Temuco"
$ws.Range("C208").Value = "La Araucanía"
$ws.Range("D208").Value = 44516
$ws.Range("E208").Value = 9
$ws.Range("F208").Value = 100112009
$ws.Range("G208").Value = "Acelga"
$ws.Range("H208").Value = "Sin especificar"
$ws.Range("I208").Value = "Primera"
$ws.Range("J208").Value = 80
$ws.Range("K208").Value = 8000
$ws.Range("L208").Value = 9000
$ws.Range("M208").Value = 8438
$ws.Range("N208").Value = "$/docena de atados (12 kilos)"
$ws.Range("O208").Value = "Provincia de Cautín"
$ws.Range("P208").Value = 703
$ws.Range("Q208").Value = 12
$ws.Range("R208").Value = "Hortaliza"
